# Continue filling in the daily progress tracker.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh the "Start" / "Tdays" inputs (L16/N16) that drive the PPD (P16)
# goal-pace calculation for the rest of the sheet.
$ws.Range("L16").Value = 46
$ws.Range("N16").Formula = "=24-7"

# Log Actual Pages written for the next several days.
$ws.Range("E16").Value = 46
$ws.Range("E17").Value = 46
$ws.Range("E18").Value = 46
$ws.Range("E19").Value = 46
$ws.Range("E20").Value = 49

# Note what section of the thesis those days were spent on.
$ws.Range("H16:H19").Value = "(Design) 2.2 data viz"

# Clear out the now-unused scratch column of running totals.
$ws.Range("J7:J30").ClearContents()

# Leave the selection where the author last left off.
$ws.Range("H19").Select()
